$wb = $excel.ActiveWorkbook
$ws0 = $wb.Worksheets.Item(1)
$ws1 = $wb.Worksheets.Add($ws0)
$ws1.Name = "Sheet1"

function Set-TextValue($ws, $cellRef, $val) {
    $helper = $ws.Range("Z1")
    $helper.Formula = "=TEXT(" + $val + ",""0"")"
    $helper.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $helper.Clear()
}

$ws1.Range("A1").Value = "TC"
$ws1.Range("B1").Value = "Customer_ID"
$ws1.Range("C1").Value = "PD"

Set-TextValue $ws1 "A2" "118448"
Set-TextValue $ws1 "B2" "17704590"
Set-TextValue $ws1 "C2" "1005"

Set-TextValue $ws1 "A3" "118450"
Set-TextValue $ws1 "B3" "17704591"
Set-TextValue $ws1 "C3" "1005"

Set-TextValue $ws1 "A4" "118451"
Set-TextValue $ws1 "B4" "17704592"
Set-TextValue $ws1 "C4" "1005"

Set-TextValue $ws1 "A5" "118452"
Set-TextValue $ws1 "B5" "17704593"
Set-TextValue $ws1 "C5" "1005"

Set-TextValue $ws1 "A6" "118452"
Set-TextValue $ws1 "B6" "17704594"
Set-TextValue $ws1 "C6" "1005"

Write-Host "done"
